# Add team record (Wins/Losses/Ties) columns to the PIT_2007 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (AC1, style "1":
# bold, bordered, centered) onto the three new header cells so they match
# the look of the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-48) gets the same team record: 68 wins, 94 losses,
# 0 ties.
$ws.Range("AD2:AD48").Value = 68
$ws.Range("AE2:AE48").Value = 94
$ws.Range("AF2:AF48").Value = 0
